$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells must stay text even though many look numeric.
# Force text format on each cell before assignment (so Excel does not
# coerce the literal into a number), then restore the default "Normal"
# style afterwards so the saved style index matches the original (no
# explicit s="..." on the cell).
$dCells = @("D2","D3","D5","D6","D9","D10","D14","D15","D16","D17","D18","D19","D20","D23","D24","D26","D29","D30","D32","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D48","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "44.056.92"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.326.41"
$ws.Range("E3").Value = "  +4.02%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "98.55"
$ws.Range("E5").Value = "  +4.50%  "
$ws.Range("D6").Value = "271.37"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.628"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("D10").Value = "45.82"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("E12").Value = "  -4.47%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "2.667.60"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").Value = "15.59"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "0.881"
$ws.Range("E16").Value = "  +7.27%  "
$ws.Range("D17").Value = "2.331.04"
$ws.Range("E17").Value = "  +4.14%  "
$ws.Range("D18").Value = "44.034.73"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "0.0000110"
$ws.Range("E19").Value = "  +4.61%  "
$ws.Range("D20").Value = "6.42"
$ws.Range("E20").Value = "  +3.78%  "
$ws.Range("E21").Value = "  +4.03%  "
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "240.46"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").Value = "9.30"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "2.54"
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").Value = "38.42"
$ws.Range("E30").Value = "  -5.00%  "
$ws.Range("D32").Value = "175.96"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.0364"
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.110"
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").Value = "4.47"
$ws.Range("E38").Value = "  +4.12%  "
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.43"
$ws.Range("E40").Value = "  +12.85%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.245"
$ws.Range("E41").Value = "  +8.67%  "
$ws.Range("D42").Value = "1.45"
$ws.Range("E42").Value = "  +25.21%  "
$ws.Range("D43").Value = "12.45"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("D44").Value = "62.79"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").Value = "9.14"
$ws.Range("E45").Value = "  +8.66%  "
$ws.Range("D46").Value = "5.35"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E47").Value = "  +4.11%  "
$ws.Range("D48").Value = "100.49"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").Value = "0.192"
$ws.Range("E50").Value = "  +16.29%  "
$ws.Range("D51").Value = "2.554.57"
$ws.Range("E51").Value = "  +4.01%  "

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

